$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 4 claim numbers (NroSiniestro) values in column B.
# Leading apostrophe forces Excel to store these as text (quote-prefixed),
# preserving leading zeros and any trailing spaces.
$ws.Range("B2").Value = "'0420194407302   "
$ws.Range("B3").Value = "'0420194407303"
$ws.Range("B4").Value = "'1220194200722"
$ws.Range("B5").Value = "'1120194100457"

# Move selection to B9 (single cell) as reflected in the saved sheet view
$ws.Range("B9").Select()
